# Apply crypto price/volume updates per the commit diff (Mon May 13 16:14:11 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.975.54"
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("D3").Value = "2.958.81"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'595.39"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'147.43"
$ws.Range("E6").Value = "  +1.46%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.957.12"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'7.25"
$ws.Range("E10").Value = "  +3.35%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  +7.07%  "
$ws.Range("D12").Value = "'0.445"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +6.17%  "
$ws.Range("D14").Value = "'33.24"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("E15").Value = "  -0.61%  "
$ws.Range("D16").Value = "3.453.01"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "62.930.50"
$ws.Range("E17").Value = "  +2.59%  "
$ws.Range("D18").Value = "'6.76"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "2.964.33"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'446.14"
$ws.Range("E20").Value = "  +2.62%  "
$ws.Range("D21").Value = "'13.50"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'0.670"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("D23").Value = "'7.12"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'81.77"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "'11.25"
$ws.Range("E25").Value = "  +2.50%  "
$ws.Range("D26").Value = "'2.15"
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "'11.90"
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D29").Value = "'7.30"
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "'2.19"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("D32").Value = "0.0₃0980"
$ws.Range("E32").Value = "  +9.92%  "
$ws.Range("D33").Value = "'26.59"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  -1.91%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "'3.15"
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("D38").Value = "'5.64"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").Value = "'2.07"
$ws.Range("E39").Value = "  +1.87%  "
$ws.Range("D40").Value = "'49.70"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").Value = "'8.54"
$ws.Range("E41").Value = "  -1.13%  "
$ws.Range("E42").Value = "  -4.49%  "
$ws.Range("D43").Value = "'0.282"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'40.32"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").Value = "2.716.05"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'134.90"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("E47").Value = "  -2.72%  "
$ws.Range("D48").Value = "'364.23"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'22.95"
$ws.Range("E51").Value = "  -4.31%  "
